$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update changed cells per row. D-column (Price) values are forced to text
# format before assignment (and the style reset afterward) so that Excel
# COM does not silently coerce numeric-looking strings (e.g. "240.00",
# "1.001") into actual numbers, which would lose their original text
# formatting/trailing zeros and change the cell type away from a string.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "30.513.32"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -1.16%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.914.07"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -1.29%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.001"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +0.03%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "240.00"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -1.55%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.001"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.04%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4791"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -2.59%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2848"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -3.20%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06704"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -2.71%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "19.50"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +1.17%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "102.91"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -2.17%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.07792"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +0.24%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.921.22"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -0.94%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.203"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -2.95%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.6707"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -4.60%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "275.19"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +0.03%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "30.523.22"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -1.12%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "1.000"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -0.01%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.000007483"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -3.32%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "12.64"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -3.61%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "5.386"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -3.52%  "

$ws.Range("B22").Value = "BinanceUSD"
$ws.Range("C22").Value = "https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "1.001"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +0.15%  "

$ws.Range("B23").Value = "Chainlink"
$ws.Range("C23").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "6.305"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -3.60%  "

$ws.Range("B24").Value = "Cosmos"
$ws.Range("C24").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "9.357"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -5.17%  "

$ws.Range("B25").Value = "Monero"
$ws.Range("C25").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "167.73"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.97%  "

$ws.Range("B26").Value = "EthereumClassic"
$ws.Range("C26").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "19.21"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -2.25%  "

$ws.Range("B27").Value = "LidoDAOToken"
$ws.Range("C27").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "2.082"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -4.03%  "

$ws.Range("B28").Value = "Toncoin"
$ws.Range("C28").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "1.383"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -0.62%  "

$ws.Range("B29").Value = "Stellar"
$ws.Range("C29").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.09976"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -4.53%  "

$ws.Range("B30").Value = "Filecoin"
$ws.Range("C30").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "4.582"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +0.39%  "

$ws.Range("B31").Value = "PancakeSwap"
$ws.Range("C31").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.515"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -2.72%  "

$ws.Range("B32").Value = "InternetComputer(DFINITY)"
$ws.Range("C32").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.254"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -2.99%  "

$ws.Range("B33").Value = "Hedera"
$ws.Range("C33").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.04739"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -3.23%  "

$ws.Range("B34").Value = "ImmutableX"
$ws.Range("C34").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.7275"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -4.36%  "

$ws.Range("B35").Value = "ARBITRUM"
$ws.Range("C35").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.114"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -3.39%  "

$ws.Range("B36").Value = "HuobiToken"
$ws.Range("C36").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.721"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -0.49%  "

$ws.Range("B37").Value = "VeChain"
$ws.Range("C37").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.01909"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -5.05%  "

$ws.Range("B38").Value = "MXToken"
$ws.Range("C38").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.625"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -1.30%  "

$ws.Range("B39").Value = "FraxShare"
$ws.Range("C39").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "6.326"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -3.22%  "

$ws.Range("B40").Value = "Aave"
$ws.Range("C40").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "73.87"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -5.28%  "

$ws.Range("B41").Value = "RenderToken"
$ws.Range("C41").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.960"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -6.43%  "

$ws.Range("B42").Value = "Quant"
$ws.Range("C42").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "106.64"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -1.13%  "

$ws.Range("B43").Value = "TrustWalletToken"
$ws.Range("C43").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.8618"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -5.73%  "

$ws.Range("B44").Value = "TheSandbox"
$ws.Range("C44").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.4264"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -4.06%  "

$ws.Range("B45").Value = "PaxDollar"
$ws.Range("C45").Value = "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.001"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +0.15%  "

$ws.Range("B46").Value = "Aptos"
$ws.Range("C46").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "7.399"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -4.07%  "

$ws.Range("B47").Value = "Maker"
$ws.Range("C47").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "950.10"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -5.02%  "

$ws.Range("B48").Value = "Algorand"
$ws.Range("C48").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.1204"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -3.57%  "

$ws.Range("B49").Value = "Elrond"
$ws.Range("C49").Value = "https://coinranking.com/coin/omwkOTglq+elrond-egld"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "34.68"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -4.09%  "

$ws.Range("B50").Value = "Cronos"
$ws.Range("C50").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.05807"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +0.64%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "8.768"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -4.98%  "
